$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.36737340994749
$ws.Range("H2").Value = 9.452302473339424
$ws.Range("I2").Value = 7.146166390990507

$ws.Range("E3").Value = 4.204847452378137
$ws.Range("H3").Value = 5.684591802882636
$ws.Range("I3").Value = 2.09143060934584

$ws.Range("E4").Value = 3.957198220786172
$ws.Range("H4").Value = 3.473176571249101
$ws.Range("I4").Value = 1.126633137245187

$ws.Range("E5").Value = 3.37502721062112
$ws.Range("H5").Value = 2.025145410352073
$ws.Range("I5").Value = 0.5911471402326958

$ws.Range("E6").Value = 2.63891650368598
$ws.Range("H6").Value = 1.101808447138428
$ws.Range("I6").Value = 0.2544740607618841

$ws.Range("E7").Value = 2.114078642529698
$ws.Range("H7").Value = 0.6433867657716956
$ws.Range("I7").Value = 0.2386662659489538

$ws.Range("E8").Value = 1.466625761326899
$ws.Range("H8").Value = 0.4102779839055163
$ws.Range("I8").Value = 0.1197929858967771
